# Applies scheduled market-data refresh values to the Leve profit sheets
# (currentAveragePrice* / LevePrice* / LeveProfit* columns, H:N) per sheet.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 53
$ws.Range("H53").Value = 516.5333000000001
$ws.Range("I53").Value = 245.70589
$ws.Range("J53").Value = 870.6923
$ws.Range("K53").Value = 245.70589
$ws.Range("L53").Value = 870.6923
$ws.Range("M53").Value = 391.29411
$ws.Range("N53").Value = -2144.6923
# Row 64
$ws.Range("H64").Value = 4447
$ws.Range("I64").Value = 4462
$ws.Range("J64").Value = 4425.5713
$ws.Range("K64").Value = 4462
$ws.Range("L64").Value = 4425.5713
$ws.Range("M64").Value = -4214
$ws.Range("N64").Value = -4921.5713
# Row 67
$ws.Range("H67").Value = 4447
$ws.Range("I67").Value = 4462
$ws.Range("J67").Value = 4425.5713
$ws.Range("K67").Value = 4462
$ws.Range("L67").Value = 4425.5713
$ws.Range("M67").Value = -3604
$ws.Range("N67").Value = -6141.5713
# Row 92
$ws.Range("H92").Value = 608.7917
$ws.Range("I92").Value = 585.4
$ws.Range("J92").Value = 647.7778
$ws.Range("K92").Value = 585.4
$ws.Range("L92").Value = 647.7778
$ws.Range("M92").Value = 662.6
$ws.Range("N92").Value = -3143.7778
# Row 129
$ws.Range("H129").Value = 869
$ws.Range("J129").Value = 894.3333
$ws.Range("L129").Value = 2682.9999
$ws.Range("N129").Value = -12682.9999
# Row 141
$ws.Range("H141").Value = 1547.5
$ws.Range("I141").Value = 1547.5
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 4642.5
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 537.5
$ws.Range("N141").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1683.21
$ws.Range("I32").Value = 1683.21
$ws.Range("K32").Value = 1683.21
$ws.Range("M32").Value = -1396.21
# Row 74
$ws.Range("H74").Value = 16131243
$ws.Range("I74").Value = 25000784
$ws.Range("J74").Value = 4803.727
$ws.Range("K74").Value = 25000784
$ws.Range("L74").Value = 4803.727
$ws.Range("M74").Value = -24999910
$ws.Range("N74").Value = -6551.727
# Row 77
$ws.Range("H77").Value = 16131243
$ws.Range("I77").Value = 25000784
$ws.Range("J77").Value = 4803.727
$ws.Range("K77").Value = 125003920
$ws.Range("L77").Value = 24018.635
$ws.Range("M77").Value = -124999552
$ws.Range("N77").Value = -32754.635
# Row 88
$ws.Range("H88").Value = 1999.1666
$ws.Range("I88").Value = 1900
$ws.Range("J88").Value = 2197.5
$ws.Range("K88").Value = 1900
$ws.Range("L88").Value = 2197.5
$ws.Range("M88").Value = -1494
$ws.Range("N88").Value = -3009.5
# Row 91
$ws.Range("H91").Value = 1999.1666
$ws.Range("I91").Value = 1900
$ws.Range("J91").Value = 2197.5
$ws.Range("K91").Value = 1900
$ws.Range("L91").Value = 2197.5
$ws.Range("M91").Value = -496
$ws.Range("N91").Value = -5005.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1346.0286
$ws.Range("I86").Value = 1221.2916
$ws.Range("J86").Value = 1618.1818
$ws.Range("K86").Value = 1221.2916
$ws.Range("L86").Value = 1618.1818
$ws.Range("M86").Value = -98.29160000000002
$ws.Range("N86").Value = -3864.1818
# Row 89
$ws.Range("H89").Value = 1346.0286
$ws.Range("I89").Value = 1221.2916
$ws.Range("J89").Value = 1618.1818
$ws.Range("K89").Value = 6106.458000000001
$ws.Range("L89").Value = 8090.909000000001
$ws.Range("M89").Value = -490.4580000000005
$ws.Range("N89").Value = -19322.909

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 2981.4285
$ws.Range("I62").Value = 2535
$ws.Range("J62").Value = 4097.5
$ws.Range("K62").Value = 2535
$ws.Range("L62").Value = 4097.5
$ws.Range("M62").Value = -1911
$ws.Range("N62").Value = -5345.5
# Row 65
$ws.Range("H65").Value = 2981.4285
$ws.Range("I65").Value = 2535
$ws.Range("J65").Value = 4097.5
$ws.Range("K65").Value = 12675
$ws.Range("L65").Value = 20487.5
$ws.Range("M65").Value = -9555
$ws.Range("N65").Value = -26727.5
# Row 74
$ws.Range("H74").Value = 13769
$ws.Range("J74").Value = 16119.333
$ws.Range("L74").Value = 16119.333
$ws.Range("N74").Value = -17867.333
# Row 77
$ws.Range("H77").Value = 13769
$ws.Range("J77").Value = 16119.333
$ws.Range("L77").Value = 48357.999
$ws.Range("N77").Value = -57093.999
# Row 134
$ws.Range("H134").Value = 1706.5883
$ws.Range("I134").Value = 1283.7931
$ws.Range("J134").Value = 4158.8
$ws.Range("K134").Value = 3851.379300000001
$ws.Range("L134").Value = 12476.4
$ws.Range("M134").Value = -1316.379300000001
$ws.Range("N134").Value = -17546.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 97
$ws.Range("H97").Value = 622.2
$ws.Range("J97").Value = 500
$ws.Range("L97").Value = 1500
$ws.Range("N97").Value = -2492
# Row 113
$ws.Range("H113").Value = 2133.9517
$ws.Range("I113").Value = 491.2143
$ws.Range("J113").Value = 3486.7942
$ws.Range("K113").Value = 1473.6429
$ws.Range("L113").Value = 10460.3826
$ws.Range("M113").Value = 696.3571000000002
$ws.Range("N113").Value = -14800.3826

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 10820.417
$ws.Range("I80").Value = 2981.875
$ws.Range("J80").Value = 26497.5
$ws.Range("K80").Value = 2981.875
$ws.Range("L80").Value = 26497.5
$ws.Range("M80").Value = -1983.875
$ws.Range("N80").Value = -28493.5
# Row 83
$ws.Range("H83").Value = 10820.417
$ws.Range("I83").Value = 2981.875
$ws.Range("J83").Value = 26497.5
$ws.Range("K83").Value = 14909.375
$ws.Range("L83").Value = 132487.5
$ws.Range("M83").Value = -9917.375
$ws.Range("N83").Value = -142471.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 1439.4667
$ws.Range("I136").Value = 1241.091
$ws.Range("J136").Value = 1985
$ws.Range("K136").Value = 3723.273
$ws.Range("L136").Value = 5955
$ws.Range("M136").Value = -1173.273
$ws.Range("N136").Value = -11055
# Row 139
$ws.Range("H139").Value = 59928.75
$ws.Range("J139").Value = 59928.75
$ws.Range("L139").Value = 59928.75
$ws.Range("N139").Value = -70208.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 5110.143
$ws.Range("I81").Value = 2000
$ws.Range("J81").Value = 5349.385
$ws.Range("K81").Value = 4000
$ws.Range("L81").Value = 10698.77
$ws.Range("M81").Value = -2939
$ws.Range("N81").Value = -12820.77
# Row 84
$ws.Range("H84").Value = 5110.143
$ws.Range("I84").Value = 2000
$ws.Range("J84").Value = 5349.385
$ws.Range("K84").Value = 20000
$ws.Range("L84").Value = 53493.85000000001
$ws.Range("M84").Value = -14696
$ws.Range("N84").Value = -64101.85000000001
# Row 122
$ws.Range("H122").Value = 1458.1364
$ws.Range("I122").Value = 1439.6471
$ws.Range("K122").Value = 4318.9413
$ws.Range("M122").Value = -1868.9413

